$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t>Coral recruitment and calcium carbonate (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t>CaCO</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t xml:space="preserve">₃) accretion are fundamental processes that help maintain coral reefs. Many reefs worldwide have experienced degradation, including a decrease in coral cover and biodiversity. Successful coral recruitment helps degraded reefs to recover, while </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t>CaCO</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t xml:space="preserve">₃ accretion by early successional benthic organisms maintains the topographic complexity of a coral reef system. It is therefore important to understand the processes that affect coral recruitment and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t>CaCO</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t xml:space="preserve">₃ accretion rates </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t>in order to</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t xml:space="preserve"> understand how coral reefs recover from disturbances.  The aim of this thesis was to determine how biophysical forcing factors affect coral recruitment, calcification and bioerosion on a pristine coral reef.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Content.InsertXML($xml)
